$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 54 so the former rows 54-64 become 56-66.
$ws.Rows.Item(54).Resize(2).Insert()

# New row 54 (Chirimoya, Especial)
$ws.Cells.Item(54, 1).Value = 5
$ws.Cells.Item(54, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(54, 3).Value = "Maule"
$ws.Cells.Item(54, 4).Value = 44522
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 7
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100107
$ws.Cells.Item(54, 8).Value = "Otros"
$ws.Cells.Item(54, 9).Value = 100107002
$ws.Cells.Item(54, 10).Value = "Chirimoya"
$ws.Cells.Item(54, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(54, 12).Value = "Especial"
$ws.Cells.Item(54, 13).Value = 60
$ws.Cells.Item(54, 14).Value = 25000
$ws.Cells.Item(54, 15).Value = 25000
$ws.Cells.Item(54, 16).Value = 25000
$ws.Cells.Item(54, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(54, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(54, 19).Value = 2500
$ws.Cells.Item(54, 20).Value = 10

# New row 55 (Chirimoya, Primera)
$ws.Cells.Item(55, 1).Value = 5
$ws.Cells.Item(55, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(55, 3).Value = "Maule"
$ws.Cells.Item(55, 4).Value = 44522
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 5).Value = 7
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100107
$ws.Cells.Item(55, 8).Value = "Otros"
$ws.Cells.Item(55, 9).Value = 100107002
$ws.Cells.Item(55, 10).Value = "Chirimoya"
$ws.Cells.Item(55, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 100
$ws.Cells.Item(55, 14).Value = 23000
$ws.Cells.Item(55, 15).Value = 23000
$ws.Cells.Item(55, 16).Value = 23000
$ws.Cells.Item(55, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(55, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(55, 19).Value = 2300
$ws.Cells.Item(55, 20).Value = 10
